$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previously generated rows/cells (A1:C1, A2, A3) so the sheet
# only ends up containing the new single-row list of data.
$ws.UsedRange.ClearContents()

# Re-create row 1 as a single row built from the new list of values,
# with a gap between the two entries (A1 and F1), as in the target sheet.
$ws.Cells.Item(1, 1).Value = "data0"
$ws.Cells.Item(1, 6).Value = "data1"
